$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Readdress" -> "Readdress 1", plus a new "Readdress N" row for each
# additional board, and a trailing yes/no prompt row.
$ws.Range("A23").Value = "Readdress 1"
$ws.Range("A24").Value = "Readdress 2"
$ws.Range("A25").Value = "Readdress 3"
$ws.Range("A26").Value = "Readdress 4"
$ws.Range("A27").Value = "Readdress 5"
$ws.Range("A28").Value = "Readdress 6"
$ws.Range("A29").Value = "Readdress 7"
$ws.Range("A30").Value = "Readdress 8"
$ws.Range("A31").Value = "Is there another board?"

# Column A needs to widen to fit the new, longer labels.
$ws.Range("A:A").ColumnWidth = 19.67

# Move the active selection to reflect where editing left off.
$ws.Range("F20").Select()
